# NPCEditor_data.xlsx bug-fix pass:
#   - append new paraphrase/question lines to several "question" (col C) cells
#   - fix a couple of typos / punctuation nits in "text" (col B) cells
#   - renumber several dialogue "ID" (col A) cells
# All edits are plain Value writes against the existing ID/text/question table
# on Sheet1 (A=ID, B=text, C=question), keyed off row number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(9, 3).Value = @'
Tell me about your family.
What is your family like?
What is your relationship with your family?
How well do you get along with your folks?
Who were your parents?
Do you like your parents?
how are your parents?
how is your family?
Are you close with your parents?
Do you like your family?
do you get along with your dad and mom?
'@

$ws.Cells.Item(17, 3).Value = @'
I want to join the Navy but I'm afraid of going into battle...
Is the military scary?
Do I have to fight in battles if I join the military?
Do I have to risk my life to serve my country?
Will I have to go into battle if I join?
Have you been in battle?
I don't want to die.
I am scared of going into battle.
the military scares me.
'@

$ws.Cells.Item(33, 3).Value = @'
What is the Navy doing to combat heavy alcohol use?
How is alcohol regulated?
How does the Navy deal with alcohol abuse?
What does the Navy do about alcohol?
'@

$ws.Cells.Item(35, 1).Value = 'clintanderson_A34_1_2'

$ws.Cells.Item(61, 1).Value = 'clintanderson_A60_2_1'

$ws.Cells.Item(97, 1).Value = 'clintanderson_A96_2_2'

$ws.Cells.Item(119, 3).Value = @'
Where can I find an internship that will give me skills for the engineering field?
Where can I find relevant internships?
Where can I find apprenticeships in engineering?
Where do I look for internships?
Where do I find engineering internships?
Where can I look for internships in engineering?
How do I get an internship?
How do I become an intern?
Are there lots of internships?
'@

$ws.Cells.Item(128, 1).Value = 'clintanderson_A127_3_1'

$ws.Cells.Item(131, 3).Value = @'
How did you first become interested in computer science?
what brought you into computer science?
did you always love programming computers?
do you remember your first introduction to computers?
how did you get introduced to computer science?
how did you choose your major?
how did you choose computer science?
when did you discover your passion?
How did you find your passion?
why did you go into computer science?
why did you join CS?
'@

$ws.Cells.Item(138, 3).Value = @'
Do you have siblings?
do you have any brothers or sisters?
how many kids in your family growing up?
are you an only child?
how is your relationship with your brother or sister?
'@

$ws.Cells.Item(160, 2).Value = 'The officers have a step system just like the enlisted sailors do. So for the officers you have O1 which is Ensign, you''ve probably heard those from the movies and things like that. For O2 that''s Lieutenant JG, Lieutenant Junior Grade. For O3 there''s a Lieutenant. O4 is Lieutenant Commander. O5 is Commander. O6 is a Captain and Captains are the ones, they''re basically really really high in their officer community. They''re the people that command the ship, command the aircraft carrier, those are Captains. After that, then you have people who have been Captains for a long time and they''re super senior people that you see on TV, those are the Admirals. So you have Rear Admiral Lower Half, Upper Half, Vice Admiral and Admiral. That''s O7, O8, O9, and O10.'

$ws.Cells.Item(165, 1).Value = 'clintanderson_A164_3_2'

$ws.Cells.Item(165, 3).Value = @'
Describe the hiring process.
what is it like getting into this line of work?
how do I get a job in the Navy?
what is recruitment?
how do I get hired?
How do I start with the Navy?
How do I enlist?
How do I join the navy?
Do you have to take tests to join the navy?
Do you have to be interviewed to join the navy?
Does the navy interview you?
Does the navy test you?
'@

$ws.Cells.Item(178, 2).Value = 'I chose to be a part of the navy yes because I wanted to further my career and have you know a good salary, all the benefits, but I also feel that you have to serve your country at sometime in your life. Now I know that sounds cliche, but that''s the way that I feel. I came from a really military family. My uncle was in the military, my dad was in the military, my brother was in the military, and me too, and although yes I could of just been a civilian, like I was on that path. It felt like, at least for me in particular, it felt like something would''ve been missing from my life, something I can''t be on my death bed and say that I never protected my country the way that I saw my family protect it. So it was a calling for me. Why would other people join the military especially if you don''t have military people in your family, you have to think. Well first you have to trust your politicians. You''d have to think that they''re going to war and they''re defending America for good reason. So the wars that we fight we hope will bring stability to the world. You might see on TV that you know children are bloodied and people are you know like running for their lives, running from terrorism, running from lots of different things, you know you''ve seen the news. So if you wanna be a part of bringing that to an end, then military is one way that you can at least try. Also defending the nation, now we''ve been safe for a long time. Yes we had 9-11 and yes there are domestic terrorists but compared to other countries in the world, you might realize that America is a pretty safe place to live. When you''re considering foreign powers trying to come in and influence us. So we have that strong military to thank for that. So if you want to be a part of that team you know, and have the stamp of veteran for your character for as long as you live, then the Navy would be a good spot to start.'

$ws.Cells.Item(182, 3).Value = @'
What is the gender mix of your workplace?
how many guys and how many ladies at work?
is there a decent mix of men and women in the service?
are there many of both sexes where you have been assigned?
are any girls in the Navy?
what is the gender ratio of the Navy?
is there a gender distribution in the military?
do you work with any girls?
have you met girls in the Navy?
can women join the Navy?
'@

$ws.Cells.Item(198, 1).Value = 'clintanderson_A197_4_1'

$ws.Cells.Item(208, 2).Value = 'STEM stands for Science Technology Engineering Mathematics, STEM. And it basically consists of what most people consider technical courses.'

$ws.Cells.Item(218, 1).Value = 'clintanderson_A217_4_2'

$ws.Cells.Item(228, 3).Value = @'
What is the most dangerous thing you've ever done on the job?
What danger have you taken on in your career?
Tell me about a time you were in danger because of your work.
Were you ever in danger due to your work responsibilities?
What is an example of danger you've been in throughout your career?
Were you ever afraid in the Navy?
What were you fearful of in the Navy?
Do accidents happen in the Navy?
were you ever scared at work?
'@

$ws.Cells.Item(229, 3).Value = @'
When is breaking the rules okay?
Can I break the rules?
Are there exceptions to rules?
Is it ever okay to break rules?
In what situation is it okay to break rules?
Can I bend the rules?
Did you ever bend the rules?
Do I have to follow orders?
'@

$ws.Cells.Item(234, 3).Value = @'
What were some of your fears entering college and your career?
Did you have fears going into college or your career?
What were you afraid of before you went to college?
Were you worried about anything when you were entering college or your career?
What concerns did you have prior to entering college and your career?
What were you fearful of before beginning college and your career?
what were you scared of going into college?
'@

$ws.Cells.Item(243, 3).Value = @'
What was the coolest thing you did in the Navy?
Do you have any sick stories about your Navy career?
Did you get to do cool stuff in the Navy?
What cool things did you do in the Navy?
Is being in the Navy fun?
Did you have fun in the navy?
What sticks out to you as an awesome Navy experience?
what was the most interesting thing you did in the navy?
'@

$ws.Cells.Item(244, 1).Value = 'clintanderson_A243_5_1'

$ws.Cells.Item(266, 1).Value = 'clintanderson_A265_5_2'

$ws.Cells.Item(280, 3).Value = @'
Can you drink alcohol when deployed?
Can I get drunk while on my ship?
Does the Navy allow you to drink?
How often are service members allowed to drink on deployment?
Is it possible to drink while serving in the Navy?
Does the Navy look down on alcohol on the ship?
Do you drink?
'@

$ws.Cells.Item(299, 1).Value = 'clintanderson_A298_6_1'

$ws.Cells.Item(324, 3).Value = @'
Are you human?
Are you alive?
Are you the same species as me?
Are you an AI?
'@

$ws.Cells.Item(333, 1).Value = 'clintanderson_A332_7_1'

$ws.Cells.Item(335, 3).Value = @'
What did you get your master's degree in?
what was your master's major?
what was your major in USC?
what did you do in USC?
what did you study in USC?
where did you go to graduate school?
did you go to graduate school?
what did you go to grad school for?
what did you study in grad school?
what is your graduate degree?
'@

$ws.Cells.Item(347, 3).Value = @'
What was your favorite class in college?
what classes did you like?
what classes stood out to you?
which course did you enjoy most?
what's the best class to take?
what was your favorite course?
what is your favorite school subject?
'@

$ws.Cells.Item(356, 1).Value = 'clintanderson_A355_7_2'

$ws.Cells.Item(361, 3).Value = @'
How has combat or the Navy changed you?
how has deployment changed you?
how has the Navy changed you?
what effect has the Navy had on you?
have you ever seen combat?
have you been in combat?
have you experienced comabt?
'@
